$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text for both rows (zh-cn / de-de columns) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "64ac46a3-a86e-4621-b834-67b9cc34f26b.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30b2d1f3a548552a6efdbc8ef6f3e5688adabea9/e2e/64ac46a3-a86e-4621-b834-67b9cc34f26b.md", "", "", "64ac46a3-a86e-4621-b834-67b9cc34f26b.md")
$wsZh.Range("J2").Value = "64ac46a3-a86e-4621-b834-67b9cc34f26b.8bc8e5dddd7413fd8fd74ac2fde7bcdf6e7ff977.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-17 08:23:53"

$wsZh.Range("I3").Value = "90a09726-e726-4c7c-b00a-d76ce1e05679.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30b2d1f3a548552a6efdbc8ef6f3e5688adabea9/e2e/90a09726-e726-4c7c-b00a-d76ce1e05679.md", "", "", "90a09726-e726-4c7c-b00a-d76ce1e05679.md")
$wsZh.Range("J3").Value = "90a09726-e726-4c7c-b00a-d76ce1e05679.bc14927c8579718e454f304eacd4bb52810c531a.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-17 08:23:53"

# --- de-de sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "64ac46a3-a86e-4621-b834-67b9cc34f26b.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30b2d1f3a548552a6efdbc8ef6f3e5688adabea9/e2e/64ac46a3-a86e-4621-b834-67b9cc34f26b.md", "", "", "64ac46a3-a86e-4621-b834-67b9cc34f26b.md")
$wsDe.Range("J2").Value = "64ac46a3-a86e-4621-b834-67b9cc34f26b.8bc8e5dddd7413fd8fd74ac2fde7bcdf6e7ff977.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-17 08:24:04"

$wsDe.Range("I3").Value = "90a09726-e726-4c7c-b00a-d76ce1e05679.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30b2d1f3a548552a6efdbc8ef6f3e5688adabea9/e2e/90a09726-e726-4c7c-b00a-d76ce1e05679.md", "", "", "90a09726-e726-4c7c-b00a-d76ce1e05679.md")
$wsDe.Range("J3").Value = "90a09726-e726-4c7c-b00a-d76ce1e05679.bc14927c8579718e454f304eacd4bb52810c531a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-17 08:24:04"
